$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 55
$ws1.Range("F3").Value = 784
$ws1.Range("F8").Value = 3895
$ws1.Range("F9").Value = 89
$ws1.Range("F10").Value = 4578
$ws1.Range("F11").Value = 499
$ws1.Range("G11").Value = 128
$ws1.Range("F12").Value = 1153

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 55
$ws4.Range("F3").Value = 784
$ws4.Range("F9").Value = 3895
$ws4.Range("F10").Value = 89
$ws4.Range("F11").Value = 4578
$ws4.Range("F12").Value = 499
$ws4.Range("G12").Value = 128
$ws4.Range("F13").Value = 1153
